$d = $word.ActiveDocument

$replacements = @(
    @("93÷6=", "52÷9="),
    @("38÷3=", "88÷3="),
    @("36÷6=", "21÷3="),
    @("30÷3=", "10÷5="),
    @("52÷3=", "89÷7="),
    @("19÷4=", "62÷8="),
    @("58÷3=", "98÷3="),
    @("14÷3=", "45÷3="),
    @("16÷8=", "56÷7="),
    @("59÷4=", "26÷3="),
    @("43÷5=", "11÷9="),
    @("33÷5=", "41÷6="),
    @("31÷7=", "69÷4="),
    @("69÷6=", "53÷8="),
    @("76÷9=", "75÷2="),
    @("42÷8=", "10÷9="),
    @("98÷7=", "89÷6="),
    @("71÷8=", "89÷7="),
    @("31÷3=", "67÷4="),
    @("93÷4=", "42÷4="),
    @("32÷4=", "83÷2="),
    @("63÷6=", "27÷6="),
    @("94÷4=", "30÷6="),
    @("43÷8=", "65÷5="),
    @("39÷4=", "66÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
